$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update Correspond Handoff Datetime (D4) and
# Correspond Handback DateTime (G4) for the cee9f3e0... row
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-02-17 02:51:42"
$wsZh.Range("G4").Value = "2016-02-17 02:52:24"

# "de-de" sheet: same update for the corresponding row
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-02-17 02:51:52"
$wsDe.Range("G4").Value = "2016-02-17 02:52:42"
